# Update "想去人数" (interest count) figures in both the "展览" and
# "全部类型" worksheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 8861
    "F3"  = 8288
    "F11" = 251
    "F14" = 5365
    "F22" = 175
    "F23" = 10
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
